$d = $word.ActiveDocument

# Locate the paragraph that starts the "country comparison" block
# (Rakousko, Slovensko, Norsko, Bulharsko, Slovinsko, Srbsko). This
# whole block, together with the two empty paragraphs that trail it
# (right up to the section break), is removed - leaving the empty
# "Normlnweb" paragraph that precedes the block as the new last
# paragraph of the document body.

$count = $d.Paragraphs.Count
$startIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Rakousko*HDP na hlavu*") {
        $startIndex = $i
        break
    }
}

if ($startIndex -gt 0) {
    $startPos = $d.Paragraphs.Item($startIndex).Range.Start
    $endPos = $d.Paragraphs.Item($count).Range.End
    $killRange = $d.Range($startPos, $endPos)
    $killRange.Delete()
}
